# Update the cryptos list worksheet with refreshed price/volume figures.
# All target cells hold text (not numeric) values in the source workbook,
# so NumberFormat is forced to "@" (Text) before each write to stop Excel
# from auto-coercing number-looking strings (e.g. "1.00") into numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.866.85"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.59%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.343.03"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.10%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.18%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "549.13"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.51%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.94"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.25%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.16%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.579"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.52%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.341.38"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.06%  "

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.08%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.52"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.54%  "

# Row 12
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.35%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.339"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.36%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.68"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.05%  "

# Row 15
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "60.780.34"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.43%  "

# Row 16
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.759.36"

# Row 17
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.15%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.338.50"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.10%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.68"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.10%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.10"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.68%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "315.76"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.41%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.61"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -3.72%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.36%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.01"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.22%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.174"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.71%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.36%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.96"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.23%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.41"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +4.35%  "

# Row 29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +9.12%  "

# Row 30
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.76"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.45%  "

# Row 31
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.41"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.88%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0₃0736"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.78%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.98"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.22%  "

# Row 34
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.37"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.33%  "

# Row 35
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "PolygonEcosystemToken"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.386"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.24%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.03"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.07%  "

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.01%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.10%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.17"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.23%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "326.95"

# Row 41
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.85%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.31"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.35%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "137.60"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.27%  "

# Row 44
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.23%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0944"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.98%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.34"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.29%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.571"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.94%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0499"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.13%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0218"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.74%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₆0222"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +6.78%  "

# Row 51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.57%  "
